# Update "gh-pages" output data snapshot (杭州-漫展信息.xlsx)
# Sheet 1: 展览 (Exhibitions)
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsShow    = $wb.Worksheets.Item(2)   # 演出
$wsLocal   = $wb.Worksheets.Item(3)   # 本地生活
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# ---- 展览 (sheet 1) ----
$wsExhibit.Range("F2").Value  = 29
$wsExhibit.Range("F3").Value  = 128
$wsExhibit.Range("F8").Value  = 231
$wsExhibit.Range("F10").Value = 10359
$wsExhibit.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202409/36k37syr1726802875327.jpeg"
$wsExhibit.Range("F19").Value = 186
$wsExhibit.Range("F21").Value = 224
$wsExhibit.Range("F22").Value = 1128
$wsExhibit.Range("F26").Value = 59
$wsExhibit.Range("F27").Value = 184
$wsExhibit.Range("F28").Value = 150
$wsExhibit.Range("F30").Value = 2876
$wsExhibit.Range("F31").Value = 968
$wsExhibit.Range("F32").Value = 701
$wsExhibit.Range("F36").Value = 877
$wsExhibit.Range("F40").Value = 1261
$wsExhibit.Range("F41").Value = 615
$wsExhibit.Range("F42").Value = 5298
$wsExhibit.Range("F45").Value = 117
$wsExhibit.Range("F46").Value = 190
$wsExhibit.Range("F48").Value = 4059

# ---- 演出 (sheet 2) ----
$wsShow.Range("F6").Value  = 4051
$wsShow.Range("F8").Value  = 55
$wsShow.Range("F14").Value = 152

# ---- 本地生活 (sheet 3) ----
$wsLocal.Range("F2").Value = 735

# ---- 全部类型 (sheet 4) ----
$wsAll.Range("F2").Value  = 735
$wsAll.Range("F4").Value  = 29
$wsAll.Range("F7").Value  = 128
$wsAll.Range("F13").Value = 233
$wsAll.Range("F15").Value = 10359
$wsAll.Range("I19").Value = "//i0.hdslb.com/bfs/openplatform/202409/36k37syr1726802875327.jpeg"
$wsAll.Range("F22").Value = 224
$wsAll.Range("F23").Value = 1128
$wsAll.Range("F26").Value = 4051
$wsAll.Range("F28").Value = 59
$wsAll.Range("F29").Value = 184
$wsAll.Range("F31").Value = 2876
$wsAll.Range("F32").Value = 968
$wsAll.Range("F33").Value = 55
$wsAll.Range("F35").Value = 701
$wsAll.Range("F44").Value = 117
$wsAll.Range("F45").Value = 190
